# Applies the "Add files via upload" commit:
#  1. Adds new math-function demo rows (19-33) + a SUMIF demo (41-52) to the
#     "Maths" worksheet.
#  2. Inserts a brand-new worksheet named "Sheet1" (a READY/DISPATCH tracker)
#     right before the "Roman table" worksheet, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Extend the "Maths" worksheet with more formula demonstrations
# ---------------------------------------------------------------------------
$maths = $wb.Worksheets.Item("Maths")

# Row 19 - ODD
$maths.Cells.Item(19, 3).Formula = "=FORMULATEXT(D19)"
$maths.Cells.Item(19, 4).Formula = "=ODD(F19)"
$maths.Cells.Item(19, 6).Value = 8

# Row 20 - PI
$maths.Cells.Item(20, 3).Formula = "=FORMULATEXT(D20)"
$maths.Cells.Item(20, 4).Formula = "=PI()"

# Row 21 - POWER
$maths.Cells.Item(21, 3).Formula = "=FORMULATEXT(D21)"
$maths.Cells.Item(21, 4).Formula = "=POWER(F21,G21)"
$maths.Cells.Item(21, 6).Value = 2
$maths.Cells.Item(21, 7).Value = 3

# Row 22 - PRODUCT
$maths.Cells.Item(22, 3).Formula = "=FORMULATEXT(D22)"
$maths.Cells.Item(22, 4).Formula = "=PRODUCT(F22,G22)"
$maths.Cells.Item(22, 6).Value = 8
$maths.Cells.Item(22, 7).Value = 2

# Row 23 - QUOTIENT
$maths.Cells.Item(23, 3).Formula = "=FORMULATEXT(D23)"
$maths.Cells.Item(23, 4).Formula = "=QUOTIENT(F23,G23)"
$maths.Cells.Item(23, 6).Value = 8
$maths.Cells.Item(23, 7).Value = 2

# Row 24 - RAND
$maths.Cells.Item(24, 3).Formula = "=FORMULATEXT(D24)"
$maths.Cells.Item(24, 4).Formula = "=RAND()"

# Row 25 - RANDBETWEEN
$maths.Cells.Item(25, 3).Formula = "=FORMULATEXT(D25)"
$maths.Cells.Item(25, 4).Formula = "=RANDBETWEEN(1,99)"

# Row 26 - ROMAN
$maths.Cells.Item(26, 3).Formula = "=FORMULATEXT(D26)"
$maths.Cells.Item(26, 4).Formula = "=ROMAN(19)"

# Row 27 - ROUND
$maths.Cells.Item(27, 3).Formula = "=FORMULATEXT(D27)"
$maths.Cells.Item(27, 4).Formula = "=ROUND(7.789,2)"

# Row 28 - ROUNDDOWN
$maths.Cells.Item(28, 3).Formula = "=FORMULATEXT(D28)"
$maths.Cells.Item(28, 4).Formula = "=ROUNDDOWN(5.456,2)"

# Row 29 - ROUNDUP
$maths.Cells.Item(29, 3).Formula = "=FORMULATEXT(D29)"
$maths.Cells.Item(29, 4).Formula = "=ROUNDUP(2.785,1)"

# Row 30 - SIGN (positive)
$maths.Cells.Item(30, 3).Formula = "=FORMULATEXT(D30)"
$maths.Cells.Item(30, 4).Formula = "=SIGN(E30)"
$maths.Cells.Item(30, 5).Value = 980238325

# Row 31 - SIGN (negative)
$maths.Cells.Item(31, 3).Formula = "=FORMULATEXT(D31)"
$maths.Cells.Item(31, 4).Formula = "=SIGN(E31)"
$maths.Cells.Item(31, 5).Value = -2375398257925

# Row 32 - SIGN (zero)
$maths.Cells.Item(32, 3).Formula = "=FORMULATEXT(D32)"
$maths.Cells.Item(32, 4).Formula = "=SIGN(E32)"
$maths.Cells.Item(32, 5).Value = 0

# Row 33 - SQRT
$maths.Cells.Item(33, 3).Formula = "=FORMULATEXT(D33)"
$maths.Cells.Item(33, 4).Formula = "=SQRT(E33)"
$maths.Cells.Item(33, 5).Value = 16

# A small SUMIF demonstration lower on the same sheet (rows 41-52)
$names = @("Simran", "Munna", "ajay", "ajay", "Simran", "Munna", "Munna", "ajay")
$amounts = @(100, 50, 100, 100, 200, 900, 100, 100)
for ($i = 0; $i -lt 8; $i++) {
    $r = 41 + $i
    $maths.Cells.Item($r, 5).Value = $names[$i]
    $maths.Cells.Item($r, 6).Value = $amounts[$i]
}
$maths.Cells.Item(49, 6).Formula = "=SUM(F41:F48)"
$maths.Cells.Item(51, 5).Value = "ajay"
$maths.Cells.Item(51, 6).Formula = "=SUMIF(E41:E48,""=ajay"",F41:F48)"
$maths.Cells.Item(52, 5).Formula = "=SUMIF(E41:E48,""<>ajay"",F41:F48)"

# Match the recorded view state after the edit: scrolled down, D49 selected
$maths.Range("D49").Select()

# ---------------------------------------------------------------------------
# 2) Insert the new "Sheet1" worksheet just before "Roman table"
# ---------------------------------------------------------------------------
$romanTable = $wb.Worksheets.Item("Roman table")
$sheet1 = $wb.Worksheets.Add($romanTable)
$sheet1.Name = "Sheet1"

$sheet1.Columns.Item(4).ColumnWidth = 20.109375
$sheet1.Columns.Item(5).ColumnWidth = 13
$sheet1.Columns.Item(8).ColumnWidth = 17.5546875

$sheet1.Cells.Item(3, 4).Value = "READY/DISPATCH"
$sheet1.Cells.Item(3, 5).Value = "DATE"
$sheet1.Cells.Item(3, 6).Value = "OUTPUT"

$sheet1.Cells.Item(4, 4).Value = "READY"
$sheet1.Cells.Item(9, 4).Value = "READY"

for ($r = 4; $r -le 11; $r++) {
    $cell = $sheet1.Cells.Item($r, 5)
    $cell.Formula = "=DATE(2023,1,RANDBETWEEN(1,60))"
    $cell.NumberFormat = "m/d/yyyy"
    $sheet1.Cells.Item($r, 6).Formula = "=IF(D" + $r + "=""READY"","""",IF(E" + $r + "<TODAY(),""Urgent Call"",""""))"
}

$sheet1.Cells.Item(2, 6).Formula = "=FORMULATEXT(F4)"

$sheet1.Range("F4").Select()

$excel.CalculateFullRebuild()
